# Update quarterly margin figures on the "LH" sheet to reflect the
# revised Revenue / Cost of Revenue figures for the affected quarters.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LH")

# Row 15: Gross Margin
$ws.Range("D15").Value = 0.2961
$ws.Range("E15").Value = 0.2618
$ws.Range("F15").Value = 0.27
$ws.Range("G15").Value = 0.2815

# Row 17: EBT margin
$ws.Range("F17").Value = 0.0502

# Row 19: Free Cash Flow Margin
$ws.Range("D19").Value = 0.1224
$ws.Range("E19").Value = 0.1033
$ws.Range("F19").Value = 0.0933
$ws.Range("G19").Value = 0.0911

# Row 27: EBITDA Margin
$ws.Range("D27").Value = 0.2116
$ws.Range("E27").Value = 0.1603
$ws.Range("F27").Value = 0.1613
$ws.Range("G27").Value = 0.1651

# Row 28: Operating Cash Flow Margin
$ws.Range("D28").Value = 0.1552
$ws.Range("E28").Value = 0.1395
$ws.Range("F28").Value = 0.128
$ws.Range("G28").Value = 0.125
